$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (column D) text value.
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr/string cell type) instead of converting
# them to numbers, and the cell Style is reset to "Normal" afterward so no
# stray cell formatting (e.g. Text number format) gets attached to the cell.
$updates = @{
    2  = "243.99"
    3  = "23.25"
    4  = "5.411"
    5  = "0.05980"
    6  = "3.466"
    7  = "6.525"
    8  = "0.8114"
    9  = "0.9130"
    10 = "0.1407"
    11 = "0.07382"
    12 = "0.03239"
    14 = "0.09349"
    15 = "3.861"
    16 = "0.001560"
    17 = "0.04669"
    18 = "0.0005941"
    19 = "0.006067"
    20 = "0.005009"
    21 = "0.0009830"
    22 = "0.00007799"
    23 = "3.612"
    40 = "0.03932"
    41 = "0.006234"
    42 = "0.1073"
    43 = "0.003000"
    44 = "0.006437"
    45 = "0.00005261"
    48 = "0.8501"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.Value = "'" + $updates[$row]
    $cell.Style = "Normal"
}
